$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.3270775788793804
$ws.Range("D2").Value = 0.7467009164709912

$ws.Range("C3").Value = 0.9522302400725522
$ws.Range("D3").Value = 0.3513258352377588

$ws.Range("C4").Value = 1.191219080787047
$ws.Range("D4").Value = 0.2462612393659342

$ws.Range("C5").Value = 0.01470612615852289
$ws.Range("D5").Value = 0.9883991846534199

$ws.Range("C6").Value = 0.4923446582549433
$ws.Range("D6").Value = 0.627353194577541

$ws.Range("C7").Value = 0.6837867455663227
$ws.Range("D7").Value = 0.5012492495519583

$ws.Range("C8").Value = -0.2375603988605006
$ws.Range("D8").Value = 0.8144204679642959

$ws.Range("C9").Value = 0.1589879954388519
$ws.Range("D9").Value = 0.8751288359282787

$ws.Range("C10").Value = -0.6513047116773216
$ws.Range("D10").Value = 0.5215956652492175

$ws.Range("C11").Value = -0.8123150088515994
$ws.Range("D11").Value = 0.4253139708228204
